$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q2").Value = 574130
$ws.Range("R2").Value = 6300238
$ws.Range("Z2").Value = $null
$ws.Range("AB2").Value = $null
